$d = $word.ActiveDocument

# The document has a paragraph holding the inline picture, followed by a
# couple of blank paragraphs and a paragraph that just contains the
# (mis-spelled / proof-errored) text "Ss", and finally the paragraph that
# carries the "_GoBack" bookmark. The edit removes the blank paragraphs and
# the "Ss" paragraph in between, leaving the picture paragraph immediately
# followed by the bookmark paragraph.

# Locate the paragraph that contains the lone "Ss" run.
$ssIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Ss") {
        $ssIndex = $i
        break
    }
}

if ($ssIndex -ge 2) {
    # Remove everything from the start of the paragraph right after the
    # picture paragraph through to (and including) the "Ss" paragraph, so
    # the picture paragraph becomes directly adjacent to the one that
    # follows (the bookmark paragraph).
    $startRange = $d.Paragraphs.Item(2).Range
    $endRange = $d.Paragraphs.Item($ssIndex).Range
    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}
